# Generate Report for Handback
# Updates the localization-status workbook after a handback:
#  - Status moves from "In Translation" to "Handed back: in sync with en-US"
#  - Per-language sheets (zh-cn / de-de) get the new "Latest Target File",
#    "Latest Handback File" and "Latest Handback DateTime" populated for both
#    tracked source files, with a hyperlink added to the new handback file.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9166b71a08c7bbf3019d1dcfcb7ba5ee12ff483/e2e/3bbb3882-6549-471c-ba45-a0eeea5147bb.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9166b71a08c7bbf3019d1dcfcb7ba5ee12ff483/e2e/7d001421-c36a-4aa3-9776-85d82d147dfb.md"
$mdName1 = "3bbb3882-6549-471c-ba45-a0eeea5147bb.md"
$mdName2 = "7d001421-c36a-4aa3-9776-85d82d147dfb.md"

# --- Overview sheet: refresh the per-language status cells ---------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$zhcn.Range("I2").Value = $mdName1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$zhcn.Range("J2").Value = "3bbb3882-6549-471c-ba45-a0eeea5147bb.7f1033938a3ead6c4e87ab481d423aa18ee0fc7e.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-04 06:27:13"

$zhcn.Range("I3").Value = $mdName2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null
$zhcn.Range("J3").Value = "7d001421-c36a-4aa3-9776-85d82d147dfb.4ab94a7116236b05a30e3fa49b1a4813f1630847.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-04 06:27:13"

# --- de-de sheet -------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

$dede.Range("I2").Value = $mdName1
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$dede.Range("J2").Value = "3bbb3882-6549-471c-ba45-a0eeea5147bb.7f1033938a3ead6c4e87ab481d423aa18ee0fc7e.de-de.xlf"
$dede.Range("K2").Value = "2016-09-04 06:27:20"

$dede.Range("G3").Value = "7d001421-c36a-4aa3-9776-85d82d147dfb.4ab94a7116236b05a30e3fa49b1a4813f1630847.de-de.xlf"
$dede.Range("I3").Value = $mdName2
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null
$dede.Range("J3").Value = "7d001421-c36a-4aa3-9776-85d82d147dfb.4ab94a7116236b05a30e3fa49b1a4813f1630847.de-de.xlf"
$dede.Range("K3").Value = "2016-09-04 06:27:20"
